$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '27.200.33'
Set-TextCell $ws.Range("E2") '  +0.66%  '

Set-TextCell $ws.Range("D3") '1.851.86'
Set-TextCell $ws.Range("E3") '  +1.16%  '

Set-TextCell $ws.Range("D4") '1.001'
Set-TextCell $ws.Range("E4") '  -0.44%  '

Set-TextCell $ws.Range("D5") '313.50'
Set-TextCell $ws.Range("E5") '  +0.39%  '

Set-TextCell $ws.Range("E6") '  -0.33%  '

Set-TextCell $ws.Range("D7") '0.4603'
Set-TextCell $ws.Range("E7") '  -0.22%  '

Set-TextCell $ws.Range("D8") '0.3707'
Set-TextCell $ws.Range("E8") '  -0.03%  '

Set-TextCell $ws.Range("D9") '0.07277'

Set-TextCell $ws.Range("E10") '  +0.75%  '

Set-TextCell $ws.Range("D11") '20.01'
Set-TextCell $ws.Range("E11") '  +0.86%  '

Set-TextCell $ws.Range("B12") 'WrappedEther'
Set-TextCell $ws.Range("C12") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws.Range("D12") '1.954.35'
Set-TextCell $ws.Range("E12") '  +6.96%  '

Set-TextCell $ws.Range("B13") 'TRON'
Set-TextCell $ws.Range("C13") 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell $ws.Range("D13") '0.07818'
Set-TextCell $ws.Range("E13") '  -1.48%  '

Set-TextCell $ws.Range("D14") '5.374'
Set-TextCell $ws.Range("E14") '  +0.62%  '

Set-TextCell $ws.Range("D15") '6.495'
Set-TextCell $ws.Range("E15") '  -1.03%  '

Set-TextCell $ws.Range("D16") '91.29'
Set-TextCell $ws.Range("E16") '  -0.31%  '

Set-TextCell $ws.Range("E17") '  -0.45%  '

Set-TextCell $ws.Range("D18") '0.000008916'
Set-TextCell $ws.Range("E18") '  +0.13%  '

Set-TextCell $ws.Range("D20") '14.70'
Set-TextCell $ws.Range("E20") '  -0.79%  '

Set-TextCell $ws.Range("D21") '27.232.29'
Set-TextCell $ws.Range("E21") '  +0.70%  '

Set-TextCell $ws.Range("D22") '5.076'
Set-TextCell $ws.Range("E22") '  -0.83%  '

Set-TextCell $ws.Range("D23") '10.49'
Set-TextCell $ws.Range("E23") '  -0.71%  '

Set-TextCell $ws.Range("D24") '2.085.34'
Set-TextCell $ws.Range("E24") '  +1.84%  '

Set-TextCell $ws.Range("D25") '1.947'
Set-TextCell $ws.Range("E25") '  +5.66%  '

Set-TextCell $ws.Range("D26") '151.51'
Set-TextCell $ws.Range("E26") '  -1.19%  '

Set-TextCell $ws.Range("D27") '18.35'
Set-TextCell $ws.Range("E27") '  -0.49%  '

Set-TextCell $ws.Range("D28") '2.063'
Set-TextCell $ws.Range("E28") '  +0.69%  '

Set-TextCell $ws.Range("D29") '115.54'
Set-TextCell $ws.Range("E29") '  -0.09%  '

Set-TextCell $ws.Range("D30") '5.040'
Set-TextCell $ws.Range("E30") '  -2.28%  '

Set-TextCell $ws.Range("D31") '0.08814'
Set-TextCell $ws.Range("E31") '  -1.09%  '

Set-TextCell $ws.Range("D32") '3.091'
Set-TextCell $ws.Range("E32") '  +4.32%  '

Set-TextCell $ws.Range("D33") '0.7611'
Set-TextCell $ws.Range("E33") '  +3.88%  '

Set-TextCell $ws.Range("E34") '  +3.20%  '

Set-TextCell $ws.Range("D35") '4.496'
Set-TextCell $ws.Range("E35") '  +1.29%  '

Set-TextCell $ws.Range("D36") '2.733'
Set-TextCell $ws.Range("E36") '  +9.96%  '

Set-TextCell $ws.Range("D37") '1.081'
Set-TextCell $ws.Range("E37") '  +0.84%  '

Set-TextCell $ws.Range("D38") '0.01941'
Set-TextCell $ws.Range("E38") '  -0.59%  '

Set-TextCell $ws.Range("D39") '0.05235'
Set-TextCell $ws.Range("E39") '  -0.09%  '

Set-TextCell $ws.Range("D40") '2.940'
Set-TextCell $ws.Range("E40") '  -0.08%  '

Set-TextCell $ws.Range("D41") '7.060'
Set-TextCell $ws.Range("E41") '  -0.84%  '

Set-TextCell $ws.Range("D42") '0.5095'
Set-TextCell $ws.Range("E42") '  -1.28%  '

Set-TextCell $ws.Range("D43") '0.1623'
Set-TextCell $ws.Range("E43") '  -0.23%  '

Set-TextCell $ws.Range("D44") '8.370'
Set-TextCell $ws.Range("E44") '  +1.92%  '

Set-TextCell $ws.Range("D45") '0.4775'
Set-TextCell $ws.Range("E45") '  -1.48%  '

Set-TextCell $ws.Range("E46") '  +0.97%  '

Set-TextCell $ws.Range("D47") '1.001'
Set-TextCell $ws.Range("E47") '  -0.39%  '

Set-TextCell $ws.Range("D48") '102.69'
Set-TextCell $ws.Range("E48") '  +0.25%  '

Set-TextCell $ws.Range("D49") '1.632'
Set-TextCell $ws.Range("E49") '  -0.17%  '

Set-TextCell $ws.Range("D50") '0.06215'
Set-TextCell $ws.Range("E50") '  +0.20%  '

Set-TextCell $ws.Range("D51") '65.67'
Set-TextCell $ws.Range("E51") '  +1.58%  '
